# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" sheet: bump the timestamp, update the
# numeric stats (Casos totales / Nuevos casos / Casos activos /
# Recuperados / Casos criticos / Muertes hoy / Muertes) for every
# country whose figures moved, and re-label rows whose rank changed
# (e.g. Chile now outranks China, Azerbaiyan now outranks Honduras, ...).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 28 de Mayo de 2020 a las 18:10"
$ws.Cells.Item(4, 2).Value = 1751970
$ws.Cells.Item(4, 3).Value = 6167
$ws.Cells.Item(4, 4).Value = 493533
$ws.Cells.Item(4, 5).Value = 1155961
$ws.Cells.Item(4, 7).Value = 369
$ws.Cells.Item(4, 8).Value = 102476
$ws.Cells.Item(8, 2).Value = 269127
$ws.Cells.Item(8, 3).Value = 1887
$ws.Cells.Item(8, 7).Value = 377
$ws.Cells.Item(8, 8).Value = 37837
$ws.Cells.Item(11, 2).Value = 182202
$ws.Cells.Item(11, 3).Value = 307
$ws.Cells.Item(11, 5).Value = 10450
$ws.Cells.Item(11, 7).Value = 19
$ws.Cells.Item(11, 8).Value = 8552
$ws.Cells.Item(12, 2).Value = 164936
$ws.Cells.Item(12, 3).Value = 6850
$ws.Cells.Item(12, 4).Value = 70102
$ws.Cells.Item(12, 5).Value = 90161
$ws.Cells.Item(12, 7).Value = 139
$ws.Cells.Item(12, 8).Value = 4673
$ws.Cells.Item(16, 2).Value = 88467
$ws.Cells.Item(16, 3).Value = 948
$ws.Cells.Item(16, 4).Value = 46766
$ws.Cells.Item(16, 5).Value = 34828
$ws.Cells.Item(16, 7).Value = 108
$ws.Cells.Item(16, 8).Value = 6873
$ws.Cells.Item(17, 1).Value = "Chile"
$ws.Cells.Item(17, 2).Value = 86943
$ws.Cells.Item(17, 3).Value = 4654
$ws.Cells.Item(17, 4).Value = 36150
$ws.Cells.Item(17, 5).Value = 49903
$ws.Cells.Item(17, 7).Value = 49
$ws.Cells.Item(17, 8).Value = 890
$ws.Cells.Item(18, 1).Value = "China"
$ws.Cells.Item(18, 2).Value = 82995
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 4).Value = 78288
$ws.Cells.Item(18, 5).Value = 73
$ws.Cells.Item(18, 8).Value = 4634
$ws.Cells.Item(38, 2).Value = 22825
$ws.Cells.Item(38, 3).Value = 352
$ws.Cells.Item(38, 5).Value = 11227
$ws.Cells.Item(38, 7).Value = 10
$ws.Cells.Item(38, 8).Value = 1038
$ws.Cells.Item(45, 2).Value = 16068
$ws.Cells.Item(45, 3).Value = 345
$ws.Cells.Item(45, 4).Value = 8952
$ws.Cells.Item(45, 5).Value = 6631
$ws.Cells.Item(45, 7).Value = 11
$ws.Cells.Item(45, 8).Value = 485
$ws.Cells.Item(47, 4).Value = 4617
$ws.Cells.Item(47, 5).Value = 8815
$ws.Cells.Item(47, 7).Value = 1
$ws.Cells.Item(47, 8).Value = 501
$ws.Cells.Item(55, 2).Value = 9134
$ws.Cells.Item(55, 3).Value = 48
$ws.Cells.Item(55, 4).Value = 6457
$ws.Cells.Item(55, 5).Value = 2359
$ws.Cells.Item(55, 7).Value = 1
$ws.Cells.Item(55, 8).Value = 318
$ws.Cells.Item(70, 1).Value = "Azerbaiyan"
$ws.Cells.Item(70, 2).Value = 4759
$ws.Cells.Item(70, 3).Value = 191
$ws.Cells.Item(70, 4).Value = 3022
$ws.Cells.Item(70, 5).Value = 1681
$ws.Cells.Item(70, 7).Value = 2
$ws.Cells.Item(70, 8).Value = 56
$ws.Cells.Item(71, 1).Value = "Honduras"
$ws.Cells.Item(71, 2).Value = 4640
$ws.Cells.Item(71, 3).Value = 239
$ws.Cells.Item(71, 4).Value = 506
$ws.Cells.Item(71, 5).Value = 3940
$ws.Cells.Item(71, 7).Value = 6
$ws.Cells.Item(71, 8).Value = 194
$ws.Cells.Item(74, 2).Value = 4008
$ws.Cells.Item(74, 3).Value = 7
$ws.Cells.Item(74, 4).Value = 3803
$ws.Cells.Item(74, 5).Value = 95
$ws.Cells.Item(82, 2).Value = 2906
$ws.Cells.Item(82, 3).Value = 3
$ws.Cells.Item(82, 5).Value = 1357
$ws.Cells.Item(82, 7).Value = 2
$ws.Cells.Item(82, 8).Value = 175
$ws.Cells.Item(93, 1).Value = "Somalia"
$ws.Cells.Item(93, 2).Value = 1828
$ws.Cells.Item(93, 3).Value = 97
$ws.Cells.Item(93, 4).Value = 310
$ws.Cells.Item(93, 5).Value = 1446
$ws.Cells.Item(93, 7).Value = 5
$ws.Cells.Item(93, 8).Value = 72
$ws.Cells.Item(94, 1).Value = "Islandia"
$ws.Cells.Item(94, 2).Value = 1805
$ws.Cells.Item(94, 4).Value = 1792
$ws.Cells.Item(94, 5).Value = 3
$ws.Cells.Item(94, 8).Value = 10
$ws.Cells.Item(119, 2).Value = 941
$ws.Cells.Item(119, 3).Value = 2
$ws.Cells.Item(119, 5).Value = 330
$ws.Cells.Item(128, 2).Value = 728
$ws.Cells.Item(128, 3).Value = 8
$ws.Cells.Item(128, 5).Value = 233
$ws.Cells.Item(138, 2).Value = 465
$ws.Cells.Item(138, 3).Value = 5
$ws.Cells.Item(138, 5).Value = 53
$ws.Cells.Item(158, 1).Value = "Malaui"
$ws.Cells.Item(158, 2).Value = 203
$ws.Cells.Item(158, 3).Value = 102
$ws.Cells.Item(158, 4).Value = 37
$ws.Cells.Item(158, 5).Value = 162
$ws.Cells.Item(158, 8).Value = 4
$ws.Cells.Item(159, 1).Value = "Martinica"
$ws.Cells.Item(159, 2).Value = 197
$ws.Cells.Item(159, 4).Value = 91
$ws.Cells.Item(159, 5).Value = 92
$ws.Cells.Item(159, 8).Value = 14
$ws.Cells.Item(160, 1).Value = "Islas Feroe"
$ws.Cells.Item(160, 2).Value = 187
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = 187
$ws.Cells.Item(160, 5).Value = 0
$ws.Cells.Item(161, 1).Value = "Mongolia"
$ws.Cells.Item(161, 3).Value = 13
$ws.Cells.Item(161, 4).Value = 43
$ws.Cells.Item(161, 5).Value = 118
$ws.Cells.Item(161, 8).Value = 0
$ws.Cells.Item(162, 1).Value = "Guadalupe"
$ws.Cells.Item(162, 2).Value = 161
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 115
$ws.Cells.Item(162, 5).Value = 32
$ws.Cells.Item(162, 8).Value = 14
$ws.Cells.Item(163, 1).Value = "Gibraltar"
$ws.Cells.Item(163, 2).Value = 158
$ws.Cells.Item(163, 3).Value = 1
$ws.Cells.Item(163, 4).Value = 147
$ws.Cells.Item(163, 5).Value = 11
$ws.Cells.Item(163, 8).Value = 0
$ws.Cells.Item(164, 1).Value = "Brunei"
$ws.Cells.Item(164, 2).Value = 141
$ws.Cells.Item(164, 4).Value = 138
$ws.Cells.Item(164, 5).Value = 1
$ws.Cells.Item(164, 8).Value = 2
$ws.Cells.Item(165, 1).Value = "Islas Caimanes"
$ws.Cells.Item(165, 2).Value = 140
$ws.Cells.Item(165, 5).Value = 72
$ws.Cells.Item(165, 8).Value = 1
$ws.Cells.Item(166, 1).Value = "Guyana"
$ws.Cells.Item(166, 4).Value = 67
$ws.Cells.Item(166, 5).Value = 61
$ws.Cells.Item(166, 8).Value = 11
$ws.Cells.Item(167, 1).Value = "Bermudas"
$ws.Cells.Item(167, 2).Value = 139
$ws.Cells.Item(167, 4).Value = 91
$ws.Cells.Item(167, 5).Value = 39
$ws.Cells.Item(167, 8).Value = 9
$ws.Cells.Item(168, 1).Value = "Zimbabue"
$ws.Cells.Item(168, 2).Value = 132
$ws.Cells.Item(168, 4).Value = 25
$ws.Cells.Item(168, 5).Value = 103
$ws.Cells.Item(168, 8).Value = 4
$ws.Cells.Item(169, 1).Value = "Camboya"
$ws.Cells.Item(169, 2).Value = 124
$ws.Cells.Item(169, 4).Value = 122
$ws.Cells.Item(169, 5).Value = 2
$ws.Cells.Item(169, 8).Value = 0
$ws.Cells.Item(170, 1).Value = "Siria"
$ws.Cells.Item(170, 2).Value = 122
$ws.Cells.Item(170, 3).Value = 1
$ws.Cells.Item(170, 4).Value = 43
$ws.Cells.Item(170, 5).Value = 75
$ws.Cells.Item(170, 8).Value = 4
$ws.Cells.Item(171, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(171, 2).Value = 116
$ws.Cells.Item(171, 4).Value = 108
$ws.Cells.Item(171, 5).Value = 0
$ws.Cells.Item(171, 8).Value = 8
$ws.Cells.Item(197, 1).Value = "Fiyi"
$ws.Cells.Item(197, 4).Value = 15
$ws.Cells.Item(197, 8).Value = 0
$ws.Cells.Item(198, 1).Value = "Curazao"
$ws.Cells.Item(198, 4).Value = 14
$ws.Cells.Item(198, 8).Value = 1
$ws.Cells.Item(199, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(200, 1).Value = "Santa Lucia"
$ws.Cells.Item(200, 4).Value = 18
$ws.Cells.Item(200, 8).Value = 0
$ws.Cells.Item(201, 1).Value = "Belice"
$ws.Cells.Item(201, 4).Value = 16
$ws.Cells.Item(201, 8).Value = 2
$ws.Cells.Item(210, 1).Value = "Seychelles"
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 8).Value = 0
$ws.Cells.Item(211, 1).Value = "Montserrat"
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 8).Value = 1
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1
$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0
$ws.Cells.Item(215, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(216, 1).Value = "San Bartolome"
